{"js": "// Odvolani.docx \u2014 \"Dopln\u011bn\u00ed 106 o skl\u00e1dky\"\n//\n// The only semantic change in this revision is that the \"Nad\u0159\u00edzen\u00fd\n// org\u00e1n\" (superior authority) conditional block switches from the\n// inline Jinja tag syntax `{% if \u2026 %}` / `{% endif %}` to the\n// docassemble \"paragraph\" tag syntax `{%p if \u2026 %}` / `{%p endif %}`\n// (a lone `p` is inserted right after the opening `{%`), matching how\n// every sibling conditional block in the template already works.\n// Everything else in the underlying XML diff (the many inserted\n// <w:proofErr w:type=\"gramStart|gramEnd\"/> tags and the run splits\n// around them) is purely Word's own grammar-checker bookkeeping and\n// carries no visible/textual change, so it is intentionally not\n// reproduced here.\n\nconst body = context.document.body;\n\n// 1) \"{% if nadrizeny %}\"  ->  \"{%p if nadrizeny %}\"\nconst ifResults = body.search(\"{% if nadrizeny %}\", { matchCase: true });\nifResults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < ifResults.items.length; i++) {\n  ifResults.items[i].insertText(\"{%p if nadrizeny %}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"{% endif %}\" (the standalone one that closes the block above,\n//    i.e. the paragraph whose entire text is just that tag) ->\n//    \"{%p endif %}\". There are other \"{% endif %}\" occurrences\n//    elsewhere in the document (unrelated blocks) that must stay as-is.\nconst endifResults = body.search(\"{% endif %}\", { matchCase: true });\nendifResults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < endifResults.items.length; i++) {\n  const range = endifResults.items[i];\n  const para = range.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  if (para.text === \"{% endif %}\") {\n    range.insertText(\"{%p endif %}\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Odvolani.docx \u2014 \"Doplneni 106 o skladky\"\n#\n# The only semantic change in this revision is that the \"Nadrizeny\n# organ\" (superior authority) conditional block switches from the\n# inline Jinja tag syntax \"{% if ... %}\" / \"{% endif %}\" to the\n# docassemble \"paragraph\" tag syntax \"{%p if ... %}\" / \"{%p endif %}\"\n# (a lone \"p\" is inserted right after the opening \"{%\"), matching how\n# every sibling conditional block in the template already works.\n# Everything else in the underlying XML diff (the many inserted\n# <w:proofErr w:type=\"gramStart|gramEnd\"/> tags and the run splits\n# around them) is purely Word's own grammar-checker bookkeeping and\n# carries no visible/textual change, so it is intentionally not\n# reproduced here.\n\n$d = $word.ActiveDocument\n\n# 1) \"{% if nadrizeny %}\" -> \"{%p if nadrizeny %}\" \u2014 this exact phrase\n#    is unique in the whole document, so a plain Find/Replace over the\n#    full content is safe.\n$rng1 = $d.Content\n$rng1.Find.Execute(\n  \"{% if nadrizeny %}\",\n  $true,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  \"{%p if nadrizeny %}\",\n  2\n) | Out-Null\n\n# 2) \"{% endif %}\" -> \"{%p endif %}\" for the paragraph that closes the\n#    block above. There are other \"{% endif %}\" tags elsewhere in the\n#    document (unrelated conditionals) that must be left untouched, so\n#    only the paragraph whose entire text is exactly \"{% endif %}\" is\n#    targeted.\nforeach ($p in $d.Paragraphs) {\n  $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($ptext -eq \"{% endif %}\") {\n    $rng2 = $p.Range\n    $rng2.Find.Execute(\n      \"{% endif %}\",\n      $true,\n      $false,\n      $false,\n      $false,\n      $false,\n      $true,\n      1,\n      $false,\n      \"{%p endif %}\",\n      2\n    ) | Out-Null\n  }\n}\n"}
